$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'running knee pads'
    'youth yoga pants'
    'boys lacrosse pants'
    'knee length baseball pants'
    'youth compression knee pad sleeve'
    'yoga pad for knees'
    'recovery leggings men'
    'hex leg sleeves basketball'
    'adult basketball'
    'knee pads boys youth'
    'wrestling leggings for men'
    'softball sliding shorts youth'
    'leggings sports'
    'basketball padded shorts'
    'best yoga knee pad'
    'mens yoga pants tight'
    'baseball tights'
    'volleyball catcher'
    'goalkeeper pads'
    'youth baseball pants knee'
    'paintball knee pads'
    'youth compression tights boys'
    'compression tights for boys basketball'
    'breathable knee pads'
    'mens spandex pants'
    'soccer pads boys'
    'hockey knee pads youth'
    'padded shorts basketball'
    'sliding shorts baseball youth'
    'running pads'
    'rodillera de basketball'
    'leg tights for men'
    'mens baseball pants black'
    'sliding shorts boys'
    'baseball pants knee'
    'adult softball pants'
    'catcher knee support'
    'cycling pads for men'
    'youth football girdle with knee pads'
    'best knee pads for basketball'
    'knee pads yoga'
    'compression pads'
    'knee pants boys'
    'basketball chart'
    'girl basketball knee pads'
    'knee pads for basketball girls'
    'knee pads men'
    'basketball compression pants boys'
    'mens compression knee'
    'knee pads by design'
    'good knee pads'
    'wrestling knee pad youth'
    'mens work pants with knee pad'
    'kneepad pants'
    'softball girls sliding shorts'
    'men yoga capri'
    'compression shorts lacrosse'
    'black knee pads for basketball'
    'mens 3/4 compression pants'
    'softball sliding shorts youth girls'
    'basketball knee pads girls'
    'sliding baseball shorts'
    'knee pads for big men'
    'sleeve knee pads basketball'
    'youth padded compression shorts basketball'
    'mens compression recovery pants'
    'arthritis friendly yoga'
    'youth compression pants for boys'
    'sports leggings for men'
    'long knee pads basketball'
    'best knee pads for yoga'
    'rodillera basketball'
    'basketball compression shorts'
    'mountain biking knee pads for men'
    'boys athletic tights youth'
    'youth boys compression tights'
    'boys compression pants football'
    'padded sliding shorts mens'
    'softball compression shorts'
    'youth boys compression leggings'
    'fitness gear floor guard'
    'boys compression tights basketball'
    'padded football pants youth'
    'football girdle with knee pads'
    'youth girls sliding shorts softball'
    'construction need pads'
    'softball sliding shorts for girls'
    'softball pants men'
    'pads for running'
    'big and tall compression pants'
    'hex pad knee sleeve'
    'adult softball'
    'compression knee sleeve men basketball'
    'volleyball hand protector'
    'yoga knee pad'
    'knee sleeve wrestling youth'
    'compression recovery pants'
    'below knee leggings'
    'basketball padded compression'
    'knee length basketball shorts for women'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
